{"js": "// Replace the three-digit-division problems with their updated values.\n// Each old value is unique within the document, so a plain text search\n// per pair is unambiguous.\nconst replacements = [\n  [\"588\u00f79=\", \"604\u00f74=\"],\n  [\"189\u00f79=\", \"814\u00f72=\"],\n  [\"317\u00f75=\", \"904\u00f79=\"],\n  [\"979\u00f76=\", \"930\u00f78=\"],\n  [\"312\u00f72=\", \"714\u00f79=\"],\n  [\"602\u00f77=\", \"256\u00f75=\"],\n  [\"846\u00f74=\", \"583\u00f76=\"],\n  [\"131\u00f74=\", \"745\u00f79=\"],\n  [\"633\u00f74=\", \"187\u00f78=\"],\n  [\"331\u00f79=\", \"105\u00f78=\"],\n  [\"556\u00f74=\", \"983\u00f74=\"],\n  [\"332\u00f73=\", \"503\u00f74=\"],\n  [\"466\u00f76=\", \"536\u00f72=\"],\n  [\"764\u00f77=\", \"856\u00f76=\"],\n  [\"558\u00f72=\", \"516\u00f79=\"],\n  [\"325\u00f76=\", \"470\u00f75=\"],\n  [\"321\u00f78=\", \"879\u00f77=\"],\n  [\"695\u00f73=\", \"710\u00f75=\"],\n  [\"991\u00f77=\", \"555\u00f72=\"],\n  [\"728\u00f78=\", \"243\u00f78=\"],\n  [\"639\u00f76=\", \"535\u00f76=\"],\n  [\"858\u00f74=\", \"502\u00f76=\"],\n  [\"260\u00f72=\", \"505\u00f72=\"],\n  [\"281\u00f78=\", \"512\u00f79=\"],\n  [\"156\u00f77=\", \"843\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit-division problems with their updated values.\n# Each old value is unique within the document, so a simple Find/Replace\n# per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"588\u00f79=\", \"604\u00f74=\"),\n    @(\"189\u00f79=\", \"814\u00f72=\"),\n    @(\"317\u00f75=\", \"904\u00f79=\"),\n    @(\"979\u00f76=\", \"930\u00f78=\"),\n    @(\"312\u00f72=\", \"714\u00f79=\"),\n    @(\"602\u00f77=\", \"256\u00f75=\"),\n    @(\"846\u00f74=\", \"583\u00f76=\"),\n    @(\"131\u00f74=\", \"745\u00f79=\"),\n    @(\"633\u00f74=\", \"187\u00f78=\"),\n    @(\"331\u00f79=\", \"105\u00f78=\"),\n    @(\"556\u00f74=\", \"983\u00f74=\"),\n    @(\"332\u00f73=\", \"503\u00f74=\"),\n    @(\"466\u00f76=\", \"536\u00f72=\"),\n    @(\"764\u00f77=\", \"856\u00f76=\"),\n    @(\"558\u00f72=\", \"516\u00f79=\"),\n    @(\"325\u00f76=\", \"470\u00f75=\"),\n    @(\"321\u00f78=\", \"879\u00f77=\"),\n    @(\"695\u00f73=\", \"710\u00f75=\"),\n    @(\"991\u00f77=\", \"555\u00f72=\"),\n    @(\"728\u00f78=\", \"243\u00f78=\"),\n    @(\"639\u00f76=\", \"535\u00f76=\"),\n    @(\"858\u00f74=\", \"502\u00f76=\"),\n    @(\"260\u00f72=\", \"505\u00f72=\"),\n    @(\"281\u00f78=\", \"512\u00f79=\"),\n    @(\"156\u00f77=\", \"843\u00f74=\")\n)\n\n$wdReplaceAll = 2\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, $wdReplaceAll)\n}\n"}
